# Extend test/utils/example.xlsx for extract_dsp_from_excel tests.
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Sheet1: add a few INDIRECT-based formulas in column D (rows 2-4).
$ws1.Range("D2").Formula = "=SUM(INDIRECT(""SINUS""))"
$ws1.Range("D3").Formula = "=SQRT(D2)"
$ws1.Range("D4").Formula = "=INDIRECT(""A1:A1"")"

# Sheet2: sum part of Sheet1 via a cross-sheet INDIRECT reference.
$ws2.Range("A1").Formula = "=SUM(INDIRECT(""Sheet1!A1:B18""))"
